# Apply updated emission figures to the "Total" and "By application" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "Total": rows 2-4, columns G-N -----------------------------
$wsTotal = $wb.Worksheets.Item("Total")

$wsTotal.Range("G2").Value = 68.48
$wsTotal.Range("H2").Value = 483.8
$wsTotal.Range("I2").Value = 952.3
$wsTotal.Range("J2").Value = 0.1633
$wsTotal.Range("K2").Value = 677.2
$wsTotal.Range("L2").Value = 1065
$wsTotal.Range("M2").Value = 0.2888
$wsTotal.Range("N2").Value = 0.1185

$wsTotal.Range("G3").Value = 90.03
$wsTotal.Range("H3").Value = 296.6
$wsTotal.Range("I3").Value = 795.7
$wsTotal.Range("J3").Value = 0.1346
$wsTotal.Range("K3").Value = 537.8
$wsTotal.Range("L3").Value = 976.7
$wsTotal.Range("M3").Value = 0.3241
$wsTotal.Range("N3").Value = 0.2275

$wsTotal.Range("G4").Value = 68.48
$wsTotal.Range("H4").Value = 288.1
$wsTotal.Range("I4").Value = 647.7
$wsTotal.Range("J4").Value = 0.1111
$wsTotal.Range("K4").Value = 392.1
$wsTotal.Range("L4").Value = 791.4
$wsTotal.Range("M4").Value = 0.3946
$wsTotal.Range("N4").Value = 0.2219

# --- Sheet "By application": rows 2-13, columns AB-AC ------------------
$wsApp = $wb.Worksheets.Item("By application")

$wsApp.Range("AB2").Value = 17.61
$wsApp.Range("AC2").Value = 0.1957

$wsApp.Range("AB3").Value = 31.89
$wsApp.Range("AC3").Value = 0.3543

$wsApp.Range("AB4").Value = 14.03
$wsApp.Range("AC4").Value = 0.1169

$wsApp.Range("AB5").Value = 16.49
$wsApp.Range("AC5").Value = 0.1374

$wsApp.Range("AB6").Value = 22.51
$wsApp.Range("AC6").Value = 0.2144

$wsApp.Range("AB7").Value = 31.14
$wsApp.Range("AC7").Value = 0.2966

$wsApp.Range("AB8").Value = 15.52
$wsApp.Range("AC8").Value = 0.115

$wsApp.Range("AB9").Value = 9.9
$wsApp.Range("AC9").Value = 0.0825

$wsApp.Range("AB10").Value = 17.61
$wsApp.Range("AC10").Value = 0.1957

$wsApp.Range("AB11").Value = 24.19
$wsApp.Range("AC11").Value = 0.2688

$wsApp.Range("AB12").Value = 9.441
$wsApp.Range("AC12").Value = 0.07867

$wsApp.Range("AB13").Value = 9.821
$wsApp.Range("AC13").Value = 0.08184
